$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.370.57"
$ws.Range("E2").Value = "  +1.24%  "
$ws.Range("D3").Value = "1.938.81"
$ws.Range("E3").Value = "  -1.04%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'243.26"
$ws.Range("E5").Value = "  +0.39%  "
$ws.Range("E6").Value = "  -1.54%  "
$ws.Range("D7").Value = "'57.66"
$ws.Range("E7").Value = "  -6.85%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("E9").Value = "  +0.40%  "
$ws.Range("D10").Value = "'55.61"
$ws.Range("E10").Value = "  -0.73%  "
$ws.Range("D11").Value = "'0.0836"
$ws.Range("E11").Value = "  +4.74%  "
$ws.Range("E12").Value = "  +1.07%  "
$ws.Range("D13").Value = "'0.820"
$ws.Range("E13").Value = "  -4.11%  "
$ws.Range("D14").Value = "'21.42"
$ws.Range("E14").Value = "  -1.58%  "
$ws.Range("D15").Value = "2.225.92"
$ws.Range("E15").Value = "  -0.90%  "
$ws.Range("E16").Value = "  -2.58%  "
$ws.Range("D17").Value = "'5.23"
$ws.Range("E17").Value = "  -3.25%  "
$ws.Range("D18").Value = "1.940.14"
$ws.Range("E18").Value = "  -1.07%  "
$ws.Range("D19").Value = "36.310.12"
$ws.Range("E19").Value = "  +1.44%  "
$ws.Range("D20").Value = "'69.75"
$ws.Range("E20").Value = "  -1.30%  "
$ws.Range("E21").Value = "  +1.32%  "
$ws.Range("D22").Value = "'229.06"
$ws.Range("E22").Value = "  -3.42%  "
$ws.Range("D23").Value = "'5.05"
$ws.Range("E23").Value = "  -2.41%  "
$ws.Range("E24").Value = "  -0.39%  "
$ws.Range("E25").Value = "  -3.71%  "
$ws.Range("E26").Value = "  +0.33%  "
$ws.Range("E27").Value = "  -5.07%  "
$ws.Range("D28").Value = "'162.01"
$ws.Range("E28").Value = "  +2.31%  "
$ws.Range("D29").Value = "'19.36"
$ws.Range("E29").Value = "  -1.56%  "
$ws.Range("E30").Value = "  -6.19%  "
$ws.Range("E31").Value = "  -1.37%  "
$ws.Range("E32").Value = "  +1.11%  "
$ws.Range("D33").Value = "'4.67"
$ws.Range("E33").Value = "  -3.96%  "
$ws.Range("E34").Value = "  +1.81%  "
$ws.Range("E35").Value = "  -2.52%  "
$ws.Range("D36").Value = "'6.20"
$ws.Range("E36").Value = "  -1.54%  "
$ws.Range("E37").Value = "  +0.22%  "
$ws.Range("D38").Value = "'1.77"
$ws.Range("E38").Value = "  -3.45%  "
$ws.Range("E39").Value = "  -7.35%  "
$ws.Range("D40").Value = "'3.04"
$ws.Range("E40").Value = "  -3.55%  "
$ws.Range("D41").Value = "'0.0975"
$ws.Range("E41").Value = "  -1.18%  "
$ws.Range("E42").Value = "  +4.94%  "
$ws.Range("D43").Value = "'1.17"
$ws.Range("E43").Value = "  -3.97%  "
$ws.Range("E44").Value = "  -1.22%  "
$ws.Range("E45").Value = "  -0.45%  "
$ws.Range("D46").Value = "1.350.96"
$ws.Range("E46").Value = "  +0.89%  "
$ws.Range("E47").Value = "  -5.13%  "
$ws.Range("D48").Value = "'87.46"
$ws.Range("E48").Value = "  -4.96%  "
$ws.Range("E49").Value = "  -4.91%  "
$ws.Range("D50").Value = "'2.81"
$ws.Range("E50").Value = "  +1.96%  "
$ws.Range("D51").Value = "'45.27"
$ws.Range("E51").Value = "  +3.54%  "
